$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = 7
$ws.Range("B43").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C43").Value = "Ñuble"
$ws.Range("D43").Value = 45002
$ws.Range("E43").Value = 16
$ws.Range("F43").Value = 100112030
$ws.Range("G43").Value = "Poroto granado"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 70
$ws.Range("K43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = 30000
$ws.Range("N43").Value = "$/saco 25 kilos"
$ws.Range("O43").Value = "Provincia de Diguillín"
$ws.Range("P43").Value = 1200
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"
